$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "a71ca460a19523728a579b6d0cd5d8f6"
$ws.Range("B13").Value = "8ad39b91ea361be147e998cfa0af63da"
$ws.Range("B16").Value = "70e7d1e7097b679f3f1fd67e7ed25765"
$ws.Range("B28").Value = "f701fc6735e33393ae9dfc09dc1f862b"
$ws.Range("B33").Value = "5a1e05adae804155f1755d426715533f"
$ws.Range("B35").Value = "ba88fc096d5b51396d06a16db6a5b223"
$ws.Range("B37").Value = "fd8aa947887b5445eec47d15d563aca6"
$ws.Range("B49").Value = "64cc81814c3cec1513bc938a1e364a49"
$ws.Range("B57").Value = "09ecd661d724a2aad45dafbb1c481fb7"
$ws.Range("B61").Value = "84e53217a5c0facdef3608d3e94e7e2d"
$ws.Range("B117").Value = "81a54580528beaa3746c0be2eb8a639f"
$ws.Range("B135").Value = "11cbe408a34939d2f06c53505a4dfbb8"
$ws.Range("B167").Value = "5f5dd08ce7f5cd9cc611da22b85780f7"
$ws.Range("B179").Value = "e366907f2668a31a23988aea7a08232d"
$ws.Range("B203").Value = "1a49785ca7a983cca90a64e4bf74ea5e"
$ws.Range("B285").Value = "3c8beecd021c7ba6ce95b9dc42b77388"
$ws.Range("B296").Value = "a9049afa33919eae41a4eefeb33ebac7"
$ws.Range("B299").Value = "ca06a29ddf84c1012ce23445464311d1"
$ws.Range("B352").Value = "12b38f88196fdcb87470b4d7f549e3d1"
$ws.Range("B360").Value = "8378c8ce3a4390b4106ae67049b24cbb"
$ws.Range("B387").Value = "a6232cc167e916c3bae255a3aa1b496e"
$ws.Range("B405").Value = "4e461a4e587a7df76456916175d9836a"
$ws.Range("B424").Value = "c3d15ba386f49a4a89cff768392ffa95"
$ws.Range("B451").Value = "13e97f56e3cd1bdb66de51c8a1c381d9"
$ws.Range("B452").Value = "598536e5b90f89c3d032d32078e9f437"
$ws.Range("B465").Value = "3b2e5448b8a6985b17c5bd20938c6886"
$ws.Range("B477").Value = "e1b8840a7130774ea1c4a2335241f85b"
$ws.Range("B483").Value = "014e24331ee73c599bdf6346c172acd0"
$ws.Range("B498").Value = "5ed290198258b61f220ff7b7a5501c71"
$ws.Range("B519").Value = "d4779c9ec8359669d2be52c94872eaaa"
$ws.Range("B552").Value = "b87b0ff9b1bd0957496b465abc3e1606"
$ws.Range("B558").Value = "6e672982fa194296a2efa95bac027c65"
$ws.Range("B563").Value = "33556c93dda0eb0c1c7a678f419a41bf"
$ws.Range("B601").Value = "1aea55cc5703b249fea06d459a96cf71"
$ws.Range("B644").Value = "58f99b83bd0f100926cc628a6bb5b9c9"
$ws.Range("B707").Value = "80e249793b7bcf8866688f4689a77af8"
$ws.Range("B712").Value = "f5c07954d5e36d9a67fc8c20c5548bcb"
$ws.Range("B740").Value = "d4374f0fa39c6f7edfbd28cca214f2b8"
$ws.Range("B765").Value = "b35c15896b2ec9c25e757ebc4578d914"
$ws.Range("B802").Value = "11e6135d92906710ca6283d07f1d1add"
$ws.Range("B839").Value = "838e687b650fda7a6da60c9e4c56a4be"
$ws.Range("B846").Value = "da70563953f6e5c1d4a1aab0bbe1d7e0"
$ws.Range("B848").Value = "661c7a2286dd8390bd5f9d2ff11d671b"
$ws.Range("B889").Value = "ec565bb99879f865a731ab258df28300"
$ws.Range("B911").Value = "cba30d7950a13a0c0967661dd8f1ded3"
$ws.Range("B962").Value = "f494afd6796008fb55083a3642c47aab"
